$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (dates) mirrors rows 617-636 exactly (same repeating date sequence),
# so copy it across to preserve literal text (avoids Excel auto-converting
# date-like strings to serial date numbers when assigned via .Value).
$ws.Range("A617:A636").Copy($ws.Range("A637:A656"))

$ws.Cells.Item(637, 3).Value = 1964.5
$ws.Cells.Item(637, 4).Value = 1806.650024414062
$ws.Cells.Item(637, 5).Value = 1050.949951171875
$ws.Cells.Item(637, 6).Value = 1766.300048828125
$ws.Cells.Item(637, 7).Value = 1749.5
$ws.Cells.Item(637, 8).Value = 38025.74987792969
$ws.Cells.Item(637, 9).Value = 0
$ws.Cells.Item(637, 10).Value = 175.697880602936

$ws.Cells.Item(638, 3).Value = 1941.25
$ws.Cells.Item(638, 4).Value = 1790.449951171875
$ws.Cells.Item(638, 5).Value = 1068.800048828125
$ws.Cells.Item(638, 6).Value = 1769.300048828125
$ws.Cells.Item(638, 7).Value = 1718.75
$ws.Cells.Item(638, 8).Value = 37895.55029296875
$ws.Cells.Item(638, 9).Value = -0.003423984678248408
$ws.Cells.Item(638, 10).Value = 175.0962937517508

$ws.Cells.Item(639, 3).Value = 1922.449951171875
$ws.Cells.Item(639, 4).Value = 1785.25
$ws.Cells.Item(639, 5).Value = 1056.199951171875
$ws.Cells.Item(639, 6).Value = 1749.699951171875
$ws.Cells.Item(639, 7).Value = 1729.550048828125
$ws.Cells.Item(639, 8).Value = 37618.74926757812
$ws.Cells.Item(639, 9).Value = -0.007304314708473397
$ws.Cells.Item(639, 10).Value = 173.8173353179008

$ws.Cells.Item(640, 3).Value = 1933.150024414062
$ws.Cells.Item(640, 4).Value = 1790.550048828125
$ws.Cells.Item(640, 5).Value = 1074.900024414062
$ws.Cells.Item(640, 6).Value = 1722.900024414062
$ws.Cells.Item(640, 7).Value = 1720.75
$ws.Cells.Item(640, 8).Value = 37761.05065917969
$ws.Cells.Item(640, 9).Value = 0.003782725220059497
$ws.Cells.Item(640, 10).Value = 174.4748385358913

$ws.Cells.Item(641, 3).Value = 1901.849975585938
$ws.Cells.Item(641, 4).Value = 1756.099975585938
$ws.Cells.Item(641, 5).Value = 1112.650024414062
$ws.Cells.Item(641, 6).Value = 1730.300048828125
$ws.Cells.Item(641, 7).Value = 1715
$ws.Cells.Item(641, 8).Value = 37713.70007324219
$ws.Cells.Item(641, 9).Value = -0.001253953084220899
$ws.Cells.Item(641, 10).Value = 174.2560552739903

$ws.Cells.Item(642, 3).Value = 1894.650024414062
$ws.Cells.Item(642, 4).Value = 1746.75
$ws.Cells.Item(642, 5).Value = 1077.550048828125
$ws.Cells.Item(642, 6).Value = 1750.400024414062
$ws.Cells.Item(642, 7).Value = 1741.199951171875
$ws.Cells.Item(642, 8).Value = 37455.65051269531
$ws.Cells.Item(642, 9).Value = -0.006842329446480399
$ws.Cells.Item(642, 10).Value = 173.0637379357615

$ws.Cells.Item(643, 3).Value = 1912.300048828125
$ws.Cells.Item(643, 4).Value = 1779.099975585938
$ws.Cells.Item(643, 5).Value = 1091
$ws.Cells.Item(643, 6).Value = 1756.349975585938
$ws.Cells.Item(643, 7).Value = 1745.150024414062
$ws.Cells.Item(643, 8).Value = 37856.65014648438
$ws.Cells.Item(643, 9).Value = 0.01070598503296977
$ws.Cells.Item(643, 10).Value = 174.9165557238516

$ws.Cells.Item(644, 3).Value = 1910.150024414062
$ws.Cells.Item(644, 4).Value = 1778.75
$ws.Cells.Item(644, 5).Value = 1077.849975585938
$ws.Cells.Item(644, 6).Value = 1789.349975585938
$ws.Cells.Item(644, 7).Value = 1782.650024414062
$ws.Cells.Item(644, 8).Value = 37910.79992675781
$ws.Cells.Item(644, 9).Value = 0.001430390170918655
$ws.Cells.Item(644, 10).Value = 175.16675464589

$ws.Cells.Item(645, 3).Value = 1950.449951171875
$ws.Cells.Item(645, 4).Value = 1807.599975585938
$ws.Cells.Item(645, 5).Value = 1083.75
$ws.Cells.Item(645, 6).Value = 1838.050048828125
$ws.Cells.Item(645, 7).Value = 1812.75
$ws.Cells.Item(645, 8).Value = 38550.34973144531
$ws.Cells.Item(645, 9).Value = 0.01686985782212682
$ws.Cells.Item(645, 10).Value = 178.1217928919295

$ws.Cells.Item(646, 3).Value = 1944.099975585938
$ws.Cells.Item(646, 4).Value = 1812.800048828125
$ws.Cells.Item(646, 5).Value = 1089.699951171875
$ws.Cells.Item(646, 6).Value = 1826.050048828125
$ws.Cells.Item(646, 7).Value = 1814.099975585938
$ws.Cells.Item(646, 8).Value = 38552.5498046875
$ws.Cells.Item(646, 9).Value = 0.00005707012407186833
$ws.Cells.Item(646, 10).Value = 178.1319583247497

$ws.Cells.Item(647, 3).Value = 1950.25
$ws.Cells.Item(647, 4).Value = 1811.849975585938
$ws.Cells.Item(647, 5).Value = 1094.650024414062
$ws.Cells.Item(647, 6).Value = 1757.849975585938
$ws.Cells.Item(647, 7).Value = 1797.199951171875
$ws.Cells.Item(647, 8).Value = 38385.89990234375
$ws.Cells.Item(647, 9).Value = -0.004322668751821118
$ws.Cells.Item(647, 10).Value = 177.3619528747986

$ws.Cells.Item(648, 3).Value = 1952.550048828125
$ws.Cells.Item(648, 4).Value = 1813.75
$ws.Cells.Item(648, 5).Value = 1080.300048828125
$ws.Cells.Item(648, 6).Value = 1741.150024414062
$ws.Cells.Item(648, 7).Value = 1848.5
$ws.Cells.Item(648, 8).Value = 38346.90075683594
$ws.Cells.Item(648, 9).Value = -0.001015975803798501
$ws.Cells.Item(648, 10).Value = 177.1817574221634

$ws.Cells.Item(649, 3).Value = 1892.150024414062
$ws.Cells.Item(649, 4).Value = 1756.5
$ws.Cells.Item(649, 5).Value = 1065.800048828125
$ws.Cells.Item(649, 6).Value = 1727.25
$ws.Cells.Item(649, 7).Value = 1805.599975585938
$ws.Cells.Item(649, 8).Value = 37454.75048828125
$ws.Cells.Item(649, 9).Value = -0.023265250931541
$ws.Cells.Item(649, 10).Value = 173.0595793752453

$ws.Cells.Item(650, 3).Value = 1894.199951171875
$ws.Cells.Item(650, 4).Value = 1736.5
$ws.Cells.Item(650, 5).Value = 1060.75
$ws.Cells.Item(650, 6).Value = 1676.449951171875
$ws.Cells.Item(650, 7).Value = 1877.449951171875
$ws.Cells.Item(650, 8).Value = 37317.94946289062
$ws.Cells.Item(650, 9).Value = -0.003652434567236724
$ws.Cells.Item(650, 10).Value = 172.4274905853437

$ws.Cells.Item(651, 3).Value = 1905.75
$ws.Cells.Item(651, 4).Value = 1760.050048828125
$ws.Cells.Item(651, 5).Value = 1114.699951171875
$ws.Cells.Item(651, 6).Value = 1662
$ws.Cells.Item(651, 7).Value = 1931.449951171875
$ws.Cells.Item(651, 8).Value = 38001.24975585938
$ws.Cells.Item(651, 9).Value = 0.01831023147850691
$ws.Cells.Item(651, 10).Value = 175.5846778512195

$ws.Cells.Item(652, 3).Value = 1896.449951171875
$ws.Cells.Item(652, 4).Value = 1752.800048828125
$ws.Cells.Item(652, 5).Value = 1106.699951171875
$ws.Cells.Item(652, 6).Value = 1692.900024414062
$ws.Cells.Item(652, 7).Value = 1920.400024414062
$ws.Cells.Item(652, 8).Value = 37915.79968261719
$ws.Cells.Item(652, 9).Value = -0.002248612184892999
$ws.Cells.Item(652, 10).Value = 175.1898560051227

$ws.Cells.Item(653, 3).Value = 1898.599975585938
$ws.Cells.Item(653, 4).Value = 1775.599975585938
$ws.Cells.Item(653, 5).Value = 1098.5
$ws.Cells.Item(653, 6).Value = 1660.900024414062
$ws.Cells.Item(653, 7).Value = 1838.75
$ws.Cells.Item(653, 8).Value = 37717.7998046875
$ws.Cells.Item(653, 9).Value = -0.00522209420840627
$ws.Cells.Item(653, 10).Value = 174.2749980727068

$ws.Cells.Item(654, 3).Value = 1895.300048828125
$ws.Cells.Item(654, 4).Value = 1782.400024414062
$ws.Cells.Item(654, 5).Value = 1088.599975585938
$ws.Cells.Item(654, 6).Value = 1654.75
$ws.Cells.Item(654, 7).Value = 1722.050048828125
$ws.Cells.Item(654, 8).Value = 37400.95031738281
$ws.Cells.Item(654, 9).Value = -0.008400529430280024
$ws.Cells.Item(654, 10).Value = 172.810995822435

$ws.Cells.Item(655, 3).Value = 1900.25
$ws.Cells.Item(655, 4).Value = 1783.849975585938
$ws.Cells.Item(655, 5).Value = 1094.949951171875
$ws.Cells.Item(655, 6).Value = 1682.449951171875
$ws.Cells.Item(655, 7).Value = 1685.900024414062
$ws.Cells.Item(655, 8).Value = 37499.49938964844
$ws.Cells.Item(655, 9).Value = 0.002634934979708856
$ws.Cells.Item(655, 10).Value = 173.2663415602059

$ws.Cells.Item(656, 3).Value = 1906.75
$ws.Cells.Item(656, 4).Value = 1808.400024414062
$ws.Cells.Item(656, 5).Value = 1121
$ws.Cells.Item(656, 6).Value = 1671
$ws.Cells.Item(656, 7).Value = 1640.800048828125
$ws.Cells.Item(656, 8).Value = 37745.10021972656
$ws.Cells.Item(656, 9).Value = 0.006549442901254356
$ws.Cells.Item(656, 10).Value = 174.4011395709637

